# Add data for 2023-11-15
# Update 2023 (column J) crime counts across the Citywide Totals sheet,
# the By Neighborhood summary sheet, and each affected neighborhood sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Andersonville")
$ws.Range("J3").Value = 3
$ws.Range("J6").Value = 9

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J6").Value = 238
$ws.Range("J7").Value = 738

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 433
$ws.Range("J3").Value = 480
$ws.Range("J4").Value = 85
$ws.Range("J6").Value = 563
$ws.Range("J7").Value = 1602

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("J6").Value = 46
$ws.Range("J7").Value = 136

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J3").Value = 9
$ws.Range("J7").Value = 738
$ws.Range("J8").Value = 1602
$ws.Range("J9").Value = 136
$ws.Range("J19").Value = 745
$ws.Range("J20").Value = 537
$ws.Range("J25").Value = 128
$ws.Range("J29").Value = 1386
$ws.Range("J31").Value = 254
$ws.Range("J33").Value = 1148
$ws.Range("J36").Value = 346
$ws.Range("J37").Value = 787
$ws.Range("J39").Value = 16
$ws.Range("J42").Value = 1101
$ws.Range("J43").Value = 218
$ws.Range("J48").Value = 291
$ws.Range("J50").Value = 152
$ws.Range("J54").Value = 489
$ws.Range("J55").Value = 390
$ws.Range("J64").Value = 169
$ws.Range("J65").Value = 637
$ws.Range("J66").Value = 78
$ws.Range("J67").Value = 956
$ws.Range("J70").Value = 37
$ws.Range("J72").Value = 99
$ws.Range("J76").Value = 374
$ws.Range("J78").Value = 299
$ws.Range("J83").Value = 510
$ws.Range("J84").Value = 212
$ws.Range("J85").Value = 1056
$ws.Range("J89").Value = 325
$ws.Range("J90").Value = 272
$ws.Range("J95").Value = 368
$ws.Range("J96").Value = 279
$ws.Range("J97").Value = 230
$ws.Range("J99").Value = 391
$ws.Range("J101").Value = 25495

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J3").Value = 213
$ws.Range("J6").Value = 288
$ws.Range("J7").Value = 745

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 149
$ws.Range("J3").Value = 183
$ws.Range("J6").Value = 152
$ws.Range("J7").Value = 537

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 6740
$ws.Range("J3").Value = 7131
$ws.Range("J4").Value = 1556
$ws.Range("J6").Value = 9509
$ws.Range("J7").Value = 25495

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J3").Value = 38
$ws.Range("J7").Value = 128

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 420
$ws.Range("J6").Value = 352
$ws.Range("J7").Value = 1386

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J2").Value = 88
$ws.Range("J7").Value = 254

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J3").Value = 382
$ws.Range("J6").Value = 398
$ws.Range("J7").Value = 1148

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J6").Value = 106
$ws.Range("J7").Value = 346

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 231
$ws.Range("J3").Value = 266
$ws.Range("J4").Value = 30
$ws.Range("J7").Value = 787

$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("J3").Value = 5
$ws.Range("J6").Value = 16

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 235
$ws.Range("J7").Value = 1101

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J2").Value = 25
$ws.Range("J6").Value = 129
$ws.Range("J7").Value = 218

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J3").Value = 53
$ws.Range("J6").Value = 143
$ws.Range("J7").Value = 291

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J6").Value = 52
$ws.Range("J7").Value = 152

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J2").Value = 120
$ws.Range("J4").Value = 35
$ws.Range("J7").Value = 489

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J3").Value = 77
$ws.Range("J7").Value = 390

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J6").Value = 59
$ws.Range("J7").Value = 169

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J3").Value = 180
$ws.Range("J7").Value = 637

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("J2").Value = 14
$ws.Range("J6").Value = 45
$ws.Range("J7").Value = 78

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 243
$ws.Range("J4").Value = 66
$ws.Range("J6").Value = 267
$ws.Range("J7").Value = 956

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("J2").Value = 16
$ws.Range("J7").Value = 37

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 99

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J3").Value = 78
$ws.Range("J7").Value = 374

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J6").Value = 93
$ws.Range("J7").Value = 299

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 150
$ws.Range("J4").Value = 20
$ws.Range("J7").Value = 510

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J6").Value = 67
$ws.Range("J7").Value = 212

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J3").Value = 369
$ws.Range("J6").Value = 307
$ws.Range("J7").Value = 1056

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J3").Value = 96
$ws.Range("J7").Value = 325

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J2").Value = 95
$ws.Range("J7").Value = 272

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J2").Value = 126
$ws.Range("J3").Value = 131
$ws.Range("J6").Value = 82
$ws.Range("J7").Value = 368

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J3").Value = 76
$ws.Range("J6").Value = 100
$ws.Range("J7").Value = 279

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J2").Value = 42
$ws.Range("J6").Value = 158
$ws.Range("J7").Value = 230

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J6").Value = 100
$ws.Range("J7").Value = 391
